$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting used by the existing data cells (style index 4, e.g. E10)
# onto the new/changed data cells in rows 10-12 first, so values keep the
# correct number format/style when set below.
$ws.Range("E9").Copy()
$ws.Range("B10:H12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 10 (Visual Basic): update B10, C10, D10; add H10
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("H10").Value = 0

# Row 11 (C): add B11, D11, F11, G11, H11
$ws.Range("B11").Value = 4
$ws.Range("D11").Value = 3
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 3

# Row 12 (PHP): add B12, H12
$ws.Range("B12").Value = 3
$ws.Range("H12").Value = 2

# Apply the same style as existing I-column cells (green fill, style index 5) to I10:I12
$ws.Range("I9").Copy()
$ws.Range("I10:I12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to B13
$ws.Range("B13").Select()
